{"js": "// Office.js (Word JavaScript API) script.\n// Rewrites the \"Biosphere Engineering\" article into the \"Exploring Music\" article,\n// keeping each paragraph's existing formatting (the replacement text is written into\n// the same paragraph range, so the run properties already on that paragraph are reused).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 0: Title\nparagraphs.items[0].getRange().insertText(\n  \"Exploring Music: Harmonizing Sounds and Emotions\",\n  Word.InsertLocation.replace\n);\n\n// 1: Author name\nparagraphs.items[1].getRange().insertText(\n  \"Isabella Clark\",\n  Word.InsertLocation.replace\n);\n\n// 2: Author email, originally split across three runs (\"sophia\" / \"oliver@biosphereengineering\" / \"org\")\n// joined by \".\" runs -- replace the whole line, which collapses cleanly to \"isabella.clark@education.com\".\nparagraphs.items[2].getRange().insertText(\n  \"isabella.clark@education.com\",\n  Word.InsertLocation.replace\n);\n\n// 3: blank separator paragraph - unchanged\n\n// 4: Main body paragraph (three \"lines\" separated by blank-line breaks, i.e. <w:br/><w:br/>,\n// represented here with the vertical-tab character \\v which Office.js maps to <w:br/>).\nconst bodyText =\n  \"Music, a universal language that transcends borders and cultures, is a captivating force in our lives.\" +\n  \" It permeates our emotions, influences our moods, and holds the power to transport us to different realms.\" +\n  \" The synergy between sounds, rhythms, and melodies weaves a rich tapestry, inviting us to delve into its intricacies and discover its mesmerizing allure.\" +\n  \" Join us on a harmonious journey as we explore the multifaceted world of music, unveiling its ability to evoke emotions, narrate stories, and connect humanity.\" +\n  \"\\v\\v\" +\n  \"In the world of music, the connection between sounds and emotions is profound.\" +\n  \" Certain melodies, harmonies, and rhythms can trigger specific emotional responses, ranging from joy and exuberance to sadness and nostalgia.\" +\n  \" This phenomenon, known as the psychology of music, has been extensively studied, revealing the intricate relationship between auditory stimuli and human emotions.\" +\n  \" Music has the uncanny ability to tap into our deepest feelings and resonate with our experiences, transporting us to a realm where emotions flow freely.\" +\n  \"\\v\\v\" +\n  \"Music is an art form that captures and reflects the human condition.\" +\n  \" It narrates stories of love, loss, triumph, and despair, weaving tales that mirror the ebb and flow of life.\" +\n  \" Through lyrics and melodies, music conveys messages that transcend words, allowing us to connect with experiences beyond our own.\" +\n  \" It captures the zeitgeist of an era, encapsulating the hopes, dreams, and fears of a generation.\" +\n  \" Music becomes a soundtrack to our lives, accompanying us through milestones and marking significant moments with its poignant melodies.\" +\n  \"\\v\\v\" +\n  \"Music transcends cultural and geographical boundaries, uniting humanity in a shared experience.\" +\n  \" It serves as a bridge between people, fostering understanding and appreciation for diverse cultures.\" +\n  \" Through shared melodies and rhythms, music creates a sense of community, bringing people together in moments of celebration, worship, and mourning.\" +\n  \" It celebrates our common humanity, reminding us that despite our differences, we are all connected by the universal language of music.\";\n\nparagraphs.items[4].getRange().insertText(bodyText, Word.InsertLocation.replace);\n\n// 5: \"Summary\" heading - unchanged\n\n// 6: Summary paragraph\nconst summaryText =\n  \"In this exploration of music, we have delved into its ability to evoke emotions, narrate stories, and connect humanity.\" +\n  \" Music is a powerful medium that communicates emotions that words cannot express, creating a tapestry of sound that resonates with our souls.\" +\n  \" It serves as a mirror to society, reflecting our joys, sorrows, and collective experiences.\" +\n  \" Moreover, music transcends cultural barriers, fostering unity and understanding among people from all walks of life.\" +\n  \" As a universal language, it brings humanity together, creating a harmonious symphony that celebrates our shared existence.\";\n\nparagraphs.items[6].getRange().insertText(summaryText, Word.InsertLocation.replace);\n\nawait context.sync();\n\n// A new empty paragraph is appended at the very end of the body.\nbody.insertParagraph(\"\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Rewrites the \"Biosphere Engineering\" article into the \"Exploring Music\" article,\n# keeping each paragraph's existing formatting (text is written back into the same\n# paragraph range, so the run properties already on that paragraph are reused).\n\n$d = $word.ActiveDocument\n\n# Helper: replace the full text of paragraph number $n (1-based, like Word COM)\n# using an explicit Range(start,end) -- setting $paragraph.Range.Text directly only\n# touches the paragraph's first run when it already holds several runs, so we rebuild\n# the Range from the paragraph's Start/End instead.\nfunction Set-ParagraphText($doc, [int]$n, [string]$text) {\n    $p = $doc.Paragraphs($n)\n    $r = $doc.Range($p.Range.Start, $p.Range.End)\n    $r.Text = $text\n}\n\n# 1: Title\nSet-ParagraphText $d 1 \"Exploring Music: Harmonizing Sounds and Emotions\"\n\n# 2: Author name\nSet-ParagraphText $d 2 \"Isabella Clark\"\n\n# 3: Author email, originally split across three runs (\"sophia\" / \"oliver@biosphereengineering\" / \"org\")\n# joined by \".\" runs -- replace the whole line, which collapses cleanly to \"isabella.clark@education.com\".\nSet-ParagraphText $d 3 \"isabella.clark@education.com\"\n\n# 4: blank separator paragraph - unchanged\n\n# 5: Main body paragraph (three \"lines\" separated by blank-line breaks, i.e. <w:br/><w:br/>,\n# represented here with Chr(11), the vertical-tab character Word maps to <w:br/>).\n$nl = [char]11\n$bodyText = (\n    \"Music, a universal language that transcends borders and cultures, is a captivating force in our lives.\" +\n    \" It permeates our emotions, influences our moods, and holds the power to transport us to different realms.\" +\n    \" The synergy between sounds, rhythms, and melodies weaves a rich tapestry, inviting us to delve into its intricacies and discover its mesmerizing allure.\" +\n    \" Join us on a harmonious journey as we explore the multifaceted world of music, unveiling its ability to evoke emotions, narrate stories, and connect humanity.\" +\n    $nl + $nl +\n    \"In the world of music, the connection between sounds and emotions is profound.\" +\n    \" Certain melodies, harmonies, and rhythms can trigger specific emotional responses, ranging from joy and exuberance to sadness and nostalgia.\" +\n    \" This phenomenon, known as the psychology of music, has been extensively studied, revealing the intricate relationship between auditory stimuli and human emotions.\" +\n    \" Music has the uncanny ability to tap into our deepest feelings and resonate with our experiences, transporting us to a realm where emotions flow freely.\" +\n    $nl + $nl +\n    \"Music is an art form that captures and reflects the human condition.\" +\n    \" It narrates stories of love, loss, triumph, and despair, weaving tales that mirror the ebb and flow of life.\" +\n    \" Through lyrics and melodies, music conveys messages that transcend words, allowing us to connect with experiences beyond our own.\" +\n    \" It captures the zeitgeist of an era, encapsulating the hopes, dreams, and fears of a generation.\" +\n    \" Music becomes a soundtrack to our lives, accompanying us through milestones and marking significant moments with its poignant melodies.\" +\n    $nl + $nl +\n    \"Music transcends cultural and geographical boundaries, uniting humanity in a shared experience.\" +\n    \" It serves as a bridge between people, fostering understanding and appreciation for diverse cultures.\" +\n    \" Through shared melodies and rhythms, music creates a sense of community, bringing people together in moments of celebration, worship, and mourning.\" +\n    \" It celebrates our common humanity, reminding us that despite our differences, we are all connected by the universal language of music.\"\n)\nSet-ParagraphText $d 5 $bodyText\n\n# 6: \"Summary\" heading - unchanged\n\n# 7: Summary paragraph\n$summaryText = (\n    \"In this exploration of music, we have delved into its ability to evoke emotions, narrate stories, and connect humanity.\" +\n    \" Music is a powerful medium that communicates emotions that words cannot express, creating a tapestry of sound that resonates with our souls.\" +\n    \" It serves as a mirror to society, reflecting our joys, sorrows, and collective experiences.\" +\n    \" Moreover, music transcends cultural barriers, fostering unity and understanding among people from all walks of life.\" +\n    \" As a universal language, it brings humanity together, creating a harmonious symphony that celebrates our shared existence.\"\n)\nSet-ParagraphText $d 7 $summaryText\n\n# A new empty paragraph is appended at the very end of the document.\n$d.Content.InsertParagraphAfter()\n"}
